# Added gift card test case
# - "Sender Details" sheet now holds the sender's own details (header text
#   renamed to senderFirst Name / senderLast Name / senderEmail /
#   senderMobile) and the mobile number is stored as a real number.
# - "Receiver Details" sheet now holds the recipient's details (header text
#   renamed to receiverFirst Name / receiverLast Name  / receiverEmail /
#   message) and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

$senderWs   = $wb.Worksheets.Item("Sender Details")
$receiverWs = $wb.Worksheets.Item("Receiver Details")

# ---- Sender Details ---------------------------------------------------
$senderWs.Range("A1").Value = "senderFirst Name"
$senderWs.Range("B1").Value = "senderLast Name"
$senderWs.Range("C1").Value = "senderEmail"
$senderWs.Range("D1").Value = "senderMobile"

$senderWs.Range("A2").Value = "Sai"
$senderWs.Range("B2").Value = "Swapna"
$senderWs.Range("C2").Value = "kattasaiswapna1@"
$senderWs.Range("D2").Value = 6304062678

# ---- Receiver Details ---------------------------------------------------
$receiverWs.Range("A1").Value = "receiverFirst Name"
$receiverWs.Range("B1").Value = "receiverLast Name "
$receiverWs.Range("C1").Value = "receiverEmail"
$receiverWs.Range("D1").Value = "message"

$receiverWs.Range("A2").Value = "Hema"
$receiverWs.Range("B2").Value = "Kumari"
$receiverWs.Range("C2").Value = "hemakumari752@gmail.com"
$receiverWs.Range("D2").Value = "Have a great day dear"

# ---- Active sheet / selection -----------------------------------------
# Previously "Sender Details" was the selected tab (selection C3); now
# "Receiver Details" becomes the selected tab with selection C11.
$receiverWs.Activate() | Out-Null
$receiverWs.Range("C11").Select() | Out-Null
